# Workbook already open; grab it and the worksheet that holds the data.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Goederenvervoer_mld_ton_eigen")

# Rename the first column's header from "Trend" to "Year" (A1).
$ws.Range("A1").Value = "Year"
